# Rename the worksheet (updates the _FilterDatabase defined name reference
# along with it, since it refers to this sheet by name).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "aux_cedis_group"

# Rename the built-in "Standard" cell style to "Normal": Excel's object
# model has no direct rename for cell styles, so recreate the style under
# the new name (which clones its formatting) and drop the old one.
$wb.Styles.Add("Normal")
$wb.Styles("Standard").Delete()

